$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = "SKU"
$ws.Range("A2").Value = 123
$ws.Range("A3").Select() | Out-Null
